$d = $word.ActiveDocument

# --- Edit 1: highlight the "consistent commits" user story green ---
$r1 = $d.Content
$r1.Find.Execute("(5 points): As a developer, I want to make good, consistent commits.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Paragraphs(1).Range.HighlightColorIndex = 4

# --- Edit 2: highlight the whole Postman user-story paragraph green ---
$r2 = $d.Content
$r2.Find.Execute("As a developer, I want to use Postman to make a POST, PUT, and both GET requests", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Paragraphs(1).Range.HighlightColorIndex = 4

# --- Edit 3: highlight the "add a new movie" user-story paragraph green, and move the
#     _GoBack bookmark from the end of the preceding paragraph into the middle of this one
#     (right after "As a film enthu") ---
$r3 = $d.Content
$r3.Find.Execute("As a film enthusiast, I want to be able to add a new movie with details, including title, genre, and director name.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para3 = $r3.Paragraphs(1)
$para3.Range.HighlightColorIndex = 4

# Locate the split point "As a film enthu" within this paragraph only (several paragraphs
# in the doc start the same way, so scope the search to this paragraph's range).
$subRange = $d.Range($para3.Range.Start, $para3.Range.End)
$subRange.Find.Execute("As a film enthu", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $subRange.End

# Re-adding the bookmark at the new collapsed range removes it from its old location
# and inserts it (splitting the run) at the new one.
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
